# docs: Add professional table formatting + intelligent page layout
#
# Replace hardcoded page breaks with style-based pagination:
#   - Heading 1 : keep-with-next + keep-lines-together
#   - Heading 2 : page-break-before + keep-with-next + keep-lines-together
#   - Heading 3 : keep-with-next + keep-lines-together
#
# This walks every paragraph in the document and, based on its paragraph
# style, applies the "Keep with next", "Keep lines together" and
# "Page break before" paragraph-format switches so that headings never get
# orphaned at the bottom of a page and top-level (Heading 2) sections always
# start on a fresh page.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal

    if ($styleName -eq "Heading 1") {
        $p.Format.KeepWithNext = $true
        $p.Format.KeepTogether = $true
    }
    elseif ($styleName -eq "Heading 2") {
        $p.Format.PageBreakBefore = $true
        $p.Format.KeepWithNext = $true
        $p.Format.KeepTogether = $true
    }
    elseif ($styleName -eq "Heading 3") {
        $p.Format.KeepWithNext = $true
        $p.Format.KeepTogether = $true
    }
}

Write-Output "Applied keep-with-next/keep-lines/page-break-before to headings"
